# Update leve-profit figures (currentAveragePrice / LevePrice / LeveProfit columns)
# across all profession sheets, per scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 119110.89
$ws.Range("I86").Value = 17999.5
$ws.Range("J86").Value = 147999.86
$ws.Range("K86").Value = 17999.5
$ws.Range("L86").Value = 147999.86
$ws.Range("M86").Value = -16876.5
$ws.Range("N86").Value = -150245.86
$ws.Range("H89").Value = 119110.89
$ws.Range("I89").Value = 17999.5
$ws.Range("J89").Value = 147999.86
$ws.Range("K89").Value = 89997.5
$ws.Range("L89").Value = 739999.2999999999
$ws.Range("M89").Value = -84381.5
$ws.Range("N89").Value = -751231.2999999999
$ws.Range("H98").Value = 1412.6471
$ws.Range("J98").Value = 6644
$ws.Range("L98").Value = 6644
$ws.Range("N98").Value = -9640
$ws.Range("H106").Value = 52664316
$ws.Range("I106").Value = 62531936
$ws.Range("J106").Value = 37000
$ws.Range("K106").Value = 62531936
$ws.Range("L106").Value = 37000
$ws.Range("M106").Value = -62531305
$ws.Range("N106").Value = -38262
$ws.Range("H122").Value = 1412.6471
$ws.Range("J122").Value = 6644
$ws.Range("L122").Value = 19932
$ws.Range("N122").Value = -24832
$ws.Range("H129").Value = 935.6842
$ws.Range("H132").Value = 2436.1724
$ws.Range("I132").Value = 967.96
$ws.Range("J132").Value = 11612.5
$ws.Range("K132").Value = 2903.88
$ws.Range("L132").Value = 34837.5
$ws.Range("M132").Value = -373.8800000000001
$ws.Range("N132").Value = -39897.5
$ws.Range("H137").Value = 5523
$ws.Range("I137").Value = 8016.8335
$ws.Range("K137").Value = 24050.5005
$ws.Range("M137").Value = -21500.5005
$ws.Range("H138").Value = 3602
$ws.Range("I138").Value = 2677.8572
$ws.Range("J138").Value = 5366.273
$ws.Range("K138").Value = 8033.571599999999
$ws.Range("L138").Value = 16098.819
$ws.Range("M138").Value = -2893.571599999999
$ws.Range("N138").Value = -26378.819
$ws.Range("H141").Value = 1863.92
$ws.Range("I141").Value = 1839.5
$ws.Range("J141").Value = 2450
$ws.Range("K141").Value = 5518.5
$ws.Range("L141").Value = 7350
$ws.Range("M141").Value = -338.5
$ws.Range("N141").Value = -17710

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 47500
$ws.Range("J24").Value = 47500
$ws.Range("L24").Value = 47500
$ws.Range("N24").Value = -48248
$ws.Range("H100").Value = 47500
$ws.Range("J100").Value = 47500
$ws.Range("L100").Value = 47500
$ws.Range("N100").Value = -49664
$ws.Range("H122").Value = 4337.6665
$ws.Range("I122").Value = 3623.4
$ws.Range("K122").Value = 10870.2
$ws.Range("M122").Value = -8420.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1171.4
$ws.Range("J64").Value = 1160.6666
$ws.Range("L64").Value = 1160.6666
$ws.Range("N64").Value = -1610.6666
$ws.Range("H67").Value = 1171.4
$ws.Range("J67").Value = 1160.6666
$ws.Range("L67").Value = 1160.6666
$ws.Range("N67").Value = -2720.6666
$ws.Range("H94").Value = 1038.25
$ws.Range("J94").Value = 719
$ws.Range("L94").Value = 719
$ws.Range("N94").Value = -1621
$ws.Range("H99").Value = 1650.2727
$ws.Range("I99").Value = 1645.3
$ws.Range("K99").Value = 1645.3
$ws.Range("M99").Value = -147.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3211.5386
$ws.Range("I58").Value = 2604.3157
$ws.Range("J58").Value = 4859.7144
$ws.Range("K58").Value = 2604.3157
$ws.Range("L58").Value = 4859.7144
$ws.Range("M58").Value = -2401.3157
$ws.Range("N58").Value = -5265.7144
$ws.Range("H62").Value = 63571.715
$ws.Range("I62").Value = 8499.25
$ws.Range("J62").Value = 137001.67
$ws.Range("K62").Value = 8499.25
$ws.Range("L62").Value = 137001.67
$ws.Range("M62").Value = -7875.25
$ws.Range("N62").Value = -138249.67
$ws.Range("H65").Value = 63571.715
$ws.Range("I65").Value = 8499.25
$ws.Range("J65").Value = 137001.67
$ws.Range("K65").Value = 42496.25
$ws.Range("L65").Value = 685008.3500000001
$ws.Range("M65").Value = -39376.25
$ws.Range("N65").Value = -691248.3500000001
$ws.Range("H136").Value = 3211.5386
$ws.Range("I136").Value = 2604.3157
$ws.Range("J136").Value = 4859.7144
$ws.Range("K136").Value = 7812.9471
$ws.Range("L136").Value = 14579.1432
$ws.Range("M136").Value = -5262.9471
$ws.Range("N136").Value = -19679.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 636.25
$ws.Range("J5").Value = 600
$ws.Range("L5").Value = 1800
$ws.Range("N5").Value = -2024
$ws.Range("H68").Value = 2089.7917
$ws.Range("I68").Value = 1721.4445
$ws.Range("J68").Value = 2310.8
$ws.Range("K68").Value = 5164.333500000001
$ws.Range("L68").Value = 6932.400000000001
$ws.Range("M68").Value = -4353.333500000001
$ws.Range("N68").Value = -8554.400000000001
$ws.Range("H71").Value = 2089.7917
$ws.Range("I71").Value = 1721.4445
$ws.Range("J71").Value = 2310.8
$ws.Range("K71").Value = 15493.0005
$ws.Range("L71").Value = 20797.2
$ws.Range("M71").Value = -11437.0005
$ws.Range("N71").Value = -28909.2
$ws.Range("H131").Value = 1441.9
$ws.Range("I131").Value = 1115
$ws.Range("J131").Value = 1448.5714
$ws.Range("K131").Value = 3345
$ws.Range("L131").Value = 4345.7142
$ws.Range("M131").Value = 1695
$ws.Range("N131").Value = -14425.7142
$ws.Range("H132").Value = 2837.8462
$ws.Range("I132").Value = 1766
$ws.Range("J132").Value = 5249.5
$ws.Range("K132").Value = 15894
$ws.Range("L132").Value = 47245.5
$ws.Range("M132").Value = -13364
$ws.Range("N132").Value = -52305.5
$ws.Range("H135").Value = 636.25
$ws.Range("J135").Value = 600
$ws.Range("L135").Value = 5400
$ws.Range("N135").Value = -10470
$ws.Range("H140").Value = 3741.0833
$ws.Range("I140").Value = 2834.125
$ws.Range("K140").Value = 8502.375
$ws.Range("M140").Value = -3322.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 3750
$ws.Range("I41").Value = 3750
$ws.Range("K41").Value = 3750
$ws.Range("M41").Value = -3395
$ws.Range("H126").Value = 4869.1113
$ws.Range("J126").Value = 4966.6665
$ws.Range("L126").Value = 14899.9995
$ws.Range("N126").Value = -19839.9995
$ws.Range("H132").Value = 2126.9167
$ws.Range("I132").Value = 1652.3
$ws.Range("K132").Value = 4956.9
$ws.Range("M132").Value = -2426.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5259900
$ws.Range("I2").Value = 5004075
$ws.Range("J2").Value = 5600999.5
$ws.Range("K2").Value = 5004075
$ws.Range("L2").Value = 5600999.5
$ws.Range("M2").Value = -5003963
$ws.Range("N2").Value = -5601223.5
$ws.Range("H16").Value = 7917.6665
$ws.Range("I16").Value = 15777.25
$ws.Range("K16").Value = 15777.25
$ws.Range("M16").Value = -15607.25
$ws.Range("H20").Value = 5083.2915
$ws.Range("I20").Value = 4636.3184
$ws.Range("K20").Value = 4636.3184
$ws.Range("M20").Value = -4410.3184
$ws.Range("H132").Value = 4200
$ws.Range("I132").Value = 4200
$ws.Range("K132").Value = 12600
$ws.Range("M132").Value = -10070

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 64830
$ws.Range("I94").Value = 74995
$ws.Range("K94").Value = 74995
$ws.Range("M94").Value = -74094
$ws.Range("H110").Value = 135000
$ws.Range("J110").Value = 135000
$ws.Range("L110").Value = 135000
$ws.Range("N110").Value = -143180
$ws.Range("H113").Value = 2654.5
$ws.Range("J113").Value = 1797
$ws.Range("L113").Value = 5391
$ws.Range("N113").Value = -9731
$ws.Range("H132").Value = 864.5
$ws.Range("I132").Value = 829.6667
$ws.Range("J132").Value = 934.1667
$ws.Range("K132").Value = 2489.0001
$ws.Range("L132").Value = 2802.5001
$ws.Range("M132").Value = 40.9998999999998
$ws.Range("N132").Value = -7862.5001
$ws.Range("H136").Value = 6529.048
$ws.Range("I136").Value = 6105.75
$ws.Range("K136").Value = 18317.25
$ws.Range("M136").Value = -15767.25
